# Daily attendance processing - normalize "Recorded By" (column G) entries
# so that "System" is listed first / the canonical-cased "System" token
# sorts ahead of the stray lowercase "system" duplicate.
#
# Observed normalizations (exact-value based, mirrors the source diff):
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "admin@admin.com, System"                 -> "System, admin@admin.com"
#   "backup@backdoor.com, system, System"     -> "backup@backdoor.com, System, system"
#
# Cells already in the canonical order (e.g. "System, dnasr281@gmail.com",
# "backup@backdoor.com, System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "admin@admin.com, System") {
        $cell.Value = "System, admin@admin.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
